# "Add files via upload" — adds a new "O-Weapon" battle-log sheet to the
# workbook, positioned right after "Group Static" (i.e. as the 10th tab),
# and selects it. This pushes "PC Cure - Magic", "PC Cure - Item" and
# "PC All Heal" one slot later in the tab order.
#
# The new sheet is created by copying an existing, structurally identical
# battle-log sheet ("PC Cure - Item") so it naturally inherits the shared
# column widths / best-fit flags used across these sheets, then its data
# rows are overwritten with the new battle entry:
#   Row 2: Player "ROBO" using the "SMG Gun" command, targeting "Skelton"
#   Row 3: Enemy  "Skelton"

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("PC Cure - Item")
$afterSheet = $wb.Worksheets.Item("Group Static")

# Worksheet.Copy(Before, After) - place the duplicate right after "Group Static"
$template.Copy([System.Reflection.Missing]::Value, $afterSheet)

$ws = $wb.Worksheets.Item($afterSheet.Index + 1)
$ws.Name = "O-Weapon"

# Wipe the copied data rows (2 onward) but keep the shared header row intact.
$ws.Rows("2:4").ClearContents()

# Row 2 - Player entry
$ws.Range("A2").Formula = "=B2"
$ws.Cells.Item(2, 2).Value = "ROBO"
$ws.Cells.Item(2, 3).Value = "Player"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 12).Value = "SMG Gun"
$ws.Cells.Item(2, 13).Value = "Skelton"

# Row 3 - Enemy entry
$ws.Range("A3").Formula = "=B3"
$ws.Cells.Item(3, 2).Value = "Skelton"
$ws.Cells.Item(3, 3).Value = "Enemy"
$ws.Cells.Item(3, 4).Value = 1

# Matches the saved selection/active cell on the new tab.
$null = $ws.Range("M3").Select()
$null = $ws.Activate()
